# fix: corrigido bug da verificação da coluna relatório gerado
#
# The "Relatório Gerado" (Report Generated) column (M) was TRUE for every
# fiscalização row except the last one (row 5), which had been left blank
# due to a bug. This sets the missing value and leaves the cursor where
# the user continued working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the missing "Relatório Gerado" flag for the last row.
[void]($ws.Range("M5").Value = $true)

# Move the selection, matching where the user continued working.
[void]($ws.Range("C13").Select())
